$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace "Activiter IDL" with "Task NUI" across the used range (columns A and I
# carry the affected labels).
$used = $ws.UsedRange
[void]$used.Replace("Activiter IDL", "Task NUI", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)

# Update the active selection to match the new state (A7 instead of J3).
[void]$ws.Range("A7").Select()
